$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# row => (D value or $null, E value)
$updates = @{
    2  = @("27.127.47", "  +0.64%  ")
    3  = @("1.679.45", "  +0.28%  ")
    4  = @($null, "  -0.03%  ")
    5  = @("215.05", "  -0.07%  ")
    6  = @($null, "  +0.12%  ")
    7  = @($null, "  -0.05%  ")
    8  = @($null, "  +1.77%  ")
    9  = @("21.36", "  +5.38%  ")
    10 = @("0.0623", "  +0.49%  ")
    11 = @("0.0886", "  -0.28%  ")
    12 = @("1.915.61", "  +0.25%  ")
    13 = @("1.686.17", "  +0.39%  ")
    14 = @($null, "  +1.49%  ")
    15 = @("0.536", "  +1.76%  ")
    16 = @("66.19", "  +0.84%  ")
    17 = @("27.111.96", "  +0.52%  ")
    18 = @("238.46", "  +1.34%  ")
    19 = @("8.13", "  +0.31%  ")
    20 = @($null, "  +1.44%  ")
    22 = @($null, "  +1.66%  ")
    23 = @("9.46", "  +3.08%  ")
    24 = @($null, "  -2.59%  ")
    25 = @("146.67", "  +0.67%  ")
    26 = @($null, "  +0.82%  ")
    27 = @($null, "  +1.56%  ")
    28 = @($null, "  +0.36%  ")
    29 = @("0.999", "  -0.23%  ")
    30 = @("0.0499", "  +0.25%  ")
    31 = @($null, "  -0.21%  ")
    32 = @("1.559.78", "  +4.95%  ")
    33 = @("3.37", "  +1.15%  ")
    34 = @("3.21", "  +1.91%  ")
    35 = @($null, "  +1.23%  ")
    36 = @($null, "  +3.01%  ")
    37 = @($null, "  -1.11%  ")
    38 = @("0.936", "  +4.23%  ")
    39 = @($null, "  +2.61%  ")
    40 = @("1.06", "  +0.58%  ")
    41 = @("69.28", "  +2.56%  ")
    42 = @($null, "  -0.02%  ")
    43 = @($null, "  -3.18%  ")
    44 = @($null, "  -1.86%  ")
    45 = @("1.824.20", "  +0.24%  ")
    46 = @("0.781", "  +0.19%  ")
    47 = @("90.77", "  +0.08%  ")
    48 = @($null, "  +3.01%  ")
    49 = @($null, "  +0.49%  ")
    50 = @("0.105", "  +2.80%  ")
    51 = @("8.13", "  +5.07%  ")
}

foreach ($row in $updates.Keys) {
    $pair = $updates[$row]
    $dVal = $pair[0]
    $eVal = $pair[1]

    if ($null -ne $dVal) {
        $dCell = $ws.Cells.Item($row, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $dVal
    }

    $eCell = $ws.Cells.Item($row, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $eVal
}
